$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(2,4,8,12,15,21,24,26,30,33,36,38,41,42,46,49,50,51,56,59,63,64,67,70,72,74,77,78,81,85,89,94,99,105,107,109,111,116,118,121,124,126,129,131,132,135,138,139,144,145,148,152,154,159,161,168,169,175,178,182,184,185,188,193)

foreach ($r in $rows) {
    $rng = $ws.Range("G" + $r + ":N" + $r)
    $rng.Value = "N/A"
}
